# refector:Excel Test get data code optimize
#
# Collapses the 3-sheet template (Sheet1/Sheet2/Sheet3) down to a single
# "Login" sheet (the old Sheet2, renamed) holding Email/Password test data,
# with a bold header row, sized columns, the saved selection/page setup
# that the author's Excel left behind, and drops the now-unused sheets.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Drop Sheet1 and Sheet3; keep Sheet2 (becomes rId1 / sheetId 2) and rename it.
$wb.Worksheets("Sheet1").Delete()
$wb.Worksheets("Sheet3").Delete()
$ws = $wb.Worksheets("Sheet2")
$ws.Name = "Login"

# Header row.
$ws.Range("A1").Value = "Email"
$ws.Range("B1").Value = "Password"
$ws.Range("A1:B1").Font.Bold = $true

# Test-data rows.
$ws.Range("A2").Value = "testdemo2@gmail.com"
$ws.Range("B2").Value = "Admin@123"
$ws.Range("A3").Value = "testdemo3@gmail.com"
$ws.Range("B3").Value = "Admin@123"
$ws.Range("A4").Value = "testdemo4@gmail.com"
$ws.Range("B4").Value = "Admin@123"

# Column sizing (Email column wide, Password column narrower).
$ws.Columns.Item(1).ColumnWidth = 37.333333333333336
$ws.Columns.Item(2).ColumnWidth = 13.833333333333334

# Leftover UI state from the author's session: selection and portrait page setup.
$ws.Range("B11").Select()
$ws.PageSetup.Orientation = 1

Write-Output "done"
